$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Add()
$ws3.Name = "case-sensitive"
$ws3.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3 = $wb.Worksheets.Item("case-sensitive")

$ws3.Range("A1").Value = "ID"
$ws3.Range("A2").Value = "a-1"
$ws3.Range("A4").Value = "c-3"
$ws3.Range("A6").Value = "A-1"
$ws3.Range("A7").Value = "B-2"
$ws3.Range("A3").Value = "b-2"
$ws3.Range("A5").Value = "d-4"
$ws3.Range("A8").Value = "C-3"
$ws3.Range("A9").Value = "D-4"
$ws3.Range("A10").Value = "1-a"
$ws3.Range("A11").Value = "2-b"
$ws3.Range("A12").Value = "3-c"
$ws3.Range("A13").Value = "4-d"
$ws3.Range("A14").Value = "1-A"
$ws3.Range("A15").Value = "2-B"
$ws3.Range("A16").Value = "3-C"
$ws3.Range("A17").Value = "4-D"
